# Preparations for CAA talk, fix Import-Function
#
# Tabelle1 ("Import"-style data sheet): fix two data values and add a new
# "add" column (F) whose cells hold the literal text "FALSE" for every
# data row, matching what the fixed import routine now writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- fix two mis-imported numbers -----------------------------------------
$ws.Range("C19").Value = -1000
$ws.Range("B20").Value = 2

# --- new column F: header + "FALSE" literal text for every data row -------
$ws.Range("F1").Value = "add"

for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    # A leading apostrophe forces the literal token "FALSE" to be stored as
    # text (a shared string) instead of being auto-typed to the Boolean
    # FALSE; ClearFormats() then drops the quote-prefix style Excel adds,
    # leaving a plain text cell.
    $cell.Value = "'FALSE"
    $cell.ClearFormats()
}

# --- view tidy-up: selection moved, no more frozen/scrolled top row -------
$ws.Range("B21").Select()
